$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the workbook-level x15ac:absPath ("last saved from" folder) is
# Excel-managed file-system metadata, not an exposed COM/document
# property, so it isn't set here.

# Update the saved selection on the sheet
$ws.Range("I22").Select()

# Fill in column J values (Today / actual remaining effort)
$ws.Range("J7").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("J13").Value = 4
$ws.Range("J14").Value = 5
$ws.Range("J15").Value = 8
$ws.Range("J16").Value = 8

# Formulas for the ideal-burndown and actual-sum rows
$ws.Range("J18").Formula = "=I18-`$F`$18/10"
$ws.Range("J19").Formula = "=SUM(J7:J16)"
